# Rename the worksheet "in" to "alamin" as part of the script refactor
# (the data rows/columns themselves are unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")
$ws.Name = "alamin"
